$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing daily data runs through row 328 (date serial 44402, 2021-07-25).
# Append 15 more daily rows (329-343) through serial 44417 (2021-08-09) —
# "aggiornamento fino a 9 agosto 2021" — matching the layout of the
# preceding rows: column A keeps the date number-format/border/bold style
# already used for A2:A328 (copied via PasteSpecial so no new style/font
# entries are created), columns B/C/D hold plain zero values.

$startRow = 329
$startSerial = 44403
$count = 15
$endRow = $startRow + $count - 1

$ws.Range("A328").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $serial = $startSerial + $i

    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}
